# Updates cryptos list (prices / volumes / a couple of row swaps)
# matching the "Updated cryptos list ... with GitHub Actions" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '27.241.94'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +0.76%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.563.84'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +0.07%  '

$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.52%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '211.00'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +1.29%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.488'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -0.29%  '

$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -0.45%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '22.21'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +0.77%  '

$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +0.17%  '

$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -0.58%  '

$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +2.03%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.787.13'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +0.07%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.570.95'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +0.49%  '

$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +0.49%  '

$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -0.45%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '27.246.83'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +0.78%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '61.79'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -0.17%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '218.15'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +0.93%  '

$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.45'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +0.94%  '

$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.0₃0702'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -0.66%  '

$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -0.46%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.14'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +0.00%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.38'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +1.73%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.95'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +0.03%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '151.45'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -1.30%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '6.63'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +0.39%  '

$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '15.02'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -0.39%  '

$ws.Range('B28').Value = 'Stellar'
$ws.Range('C28').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.107'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +0.96%  '

$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -0.41%  '

$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +2.12%  '

$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -0.54%  '

$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +0.32%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.458.54'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +2.23%  '

$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +0.25%  '

$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +4.94%  '

$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +1.11%  '

$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +0.33%  '

$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -0.26%  '

$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +1.10%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '5.86'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -0.66%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.815'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +0.65%  '

$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -0.34%  '

$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +1.21%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.976'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -2.50%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '64.37'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -0.51%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.76'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +0.42%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.701.07'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -0.07%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '85.91'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -1.31%  '

$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -0.06%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0525'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +1.20%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0947'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -1.29%  '
